$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "2023-12-05 22:14:06"
$ws.Range("B6").Value = 0.0004
